# User_Schema_Rights_Definition.xlsx
# Add config for Dataprocessor_Submodules_Table_Description.xlsx so that
# submodule tables get created (rows 61-67 on "rights_and_functions"),
# plus the matching cell comment on K64 and updated selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

# ---------------------------------------------------------------------
# Row 61 - TABLE_DESCRIPTION / SCRIPTNAME block (like row 46)
# ---------------------------------------------------------------------
$ws.Range("A61").Value = "./R-dataprocessor/submodules/Dataprocessor_Submodules_Table_Description.xlsx[frontend_table_description] "
$ws.Range("B61").Value = "330_cre_table_datap_submodules_dataproc_in.sql"
$ws.Range("C61").ClearFormats()
$ws.Range("C61").Value = "template_cre_table.sql"
$ws.Range("D61").Value = "db2dataprocessor_user"
$ws.Range("E61").Value = "db2dataprocessor_in"
$ws.Range("I61").Value = "INSERT, DELETE, UPDATE, SELECT"
$ws.Range("J61").Value = "db2dataprocessor_user"

# ---------------------------------------------------------------------
# Row 62 - RIGHTS continuation (like row 47)
# ---------------------------------------------------------------------
$ws.Range("I62").Value = "INSERT, DELETE, UPDATE, SELECT"
$ws.Range("J62").Value = "db_user"
# C62 no longer carries the old empty formatted cell from this block
$ws.Range("C62").ClearContents()
$ws.Range("C62").ClearFormats()

# ---------------------------------------------------------------------
# Row 63 - RIGHTS continuation (like row 48)
# ---------------------------------------------------------------------
$ws.Range("I63").Value = "SELECT"
$ws.Range("J63").Value = "db_log_user"

# ---------------------------------------------------------------------
# Row 64 - db_log table block (like row 52/54)
# ---------------------------------------------------------------------
$ws.Range("B64").Value = "331_cre_table_datap_submodules_log.sql"
$ws.Range("C64").Value = "template_cre_table.sql"
$ws.Range("D64").Value = "db_log_user"
$ws.Range("E64").Value = "db_log"
$ws.Range("F64").Value = "INT_ID"
$ws.Range("I64").Value = "INSERT, DELETE, UPDATE, SELECT"
$ws.Range("J64").Value = "db_log_user"
$ws.Range("K64").Value = "332_db_submodules_dp_in_to_db_log.sql"
$ws.Range("L64").Value = "template_copy_function.sql"
$ws.Range("M64").Value = "copy_submodules_dp_in_to_db_log"
$ws.Range("N64").Value = "db2dataprocessor_in"

# Comment on K64, matching the existing "Wie 30 und 31" notes on K52/K54
$ws.Range("K64").AddComment("Autor:" + [char]10 + "Wie 30 und 31")

# ---------------------------------------------------------------------
# Row 65 - RIGHTS continuation (like row 53/55)
# ---------------------------------------------------------------------
$ws.Range("J65").Value = "db_user"

# ---------------------------------------------------------------------
# Row 66 - last_import view block (like row 59)
# ---------------------------------------------------------------------
$ws.Range("B66").Value = "334_cre_view_dataproc_submodules_last_import.sql"
$ws.Range("C66").Value = "template_cre_view_last_import.sql"
$ws.Range("D66").Value = "db2dataprocessor_user"
$ws.Range("E66").Value = "db2dataprocessor_out"
$ws.Range("G66").Value = "v_"
$ws.Range("H66").Value = "_last_import"
$ws.Range("I66").Value = "SELECT"
$ws.Range("J66").Value = "db2dataprocessor_user"
$ws.Range("N66").Value = "db_log"

# ---------------------------------------------------------------------
# Row 67 - all view block (like row 60)
# ---------------------------------------------------------------------
$ws.Range("B67").Value = "335_cre_view_dataproc_submodules_all.sql"
$ws.Range("C67").Value = "template_cre_view_all.sql"
$ws.Range("D67").Value = "db2dataprocessor_user"
$ws.Range("E67").Value = "db2dataprocessor_out"
$ws.Range("G67").Value = "v_"
$ws.Range("I67").Value = "SELECT"
$ws.Range("J67").Value = "db2dataprocessor_user"
$ws.Range("N67").Value = "db_log"

# ---------------------------------------------------------------------
# View / selection - matches the author scrolling to the new block
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A64").Select()

Write-Host "Applied Dataprocessor_Submodules_Table_Description.xlsx config (rows 61-67)"
